$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I5").Value = -0.4541269873383775
$ws.Range("J5").Value = 0.4506361868609465
$ws.Range("K5").Value = 0.07358145160110161
$ws.Range("L5").Value = 2.480718658867192
